$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "314.88")
# but the sheet stores them as plain text (inline strings), so force the
# column to Text format before writing the new quotes, then drop the
# explicit format again so the cells keep their original (default) style.
$dRange = $ws.Range("D2:D51")
$originalFormat = $dRange.NumberFormat
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.329.13'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = '1.857.68'
$ws.Range("E3").Value = '  +1.60%  '
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").Value = '314.88'
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("D7").Value = '0.4609'
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").Value = '0.3713'
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("D9").Value = '0.07301'
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("D10").Value = '0.8911'
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("D11").Value = '20.11'
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("D12").Value = '0.07822'
$ws.Range("E12").Value = '  -1.45%  '
$ws.Range("D13").Value = '1.906.43'
$ws.Range("E13").Value = '  +5.17%  '
$ws.Range("D14").Value = '5.398'
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("D15").Value = '6.541'
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").Value = '91.61'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("D18").Value = '0.000008929'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").Value = '14.78'
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").Value = '27.354.36'
$ws.Range("E21").Value = '  +1.56%  '
$ws.Range("D22").Value = '5.128'
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("D24").Value = '2.072.89'
$ws.Range("E24").Value = '  -3.69%  '
$ws.Range("D25").Value = '1.927'
$ws.Range("E25").Value = '  +4.39%  '
$ws.Range("D26").Value = '152.30'
$ws.Range("E26").Value = '  -0.44%  '
$ws.Range("D27").Value = '18.48'
$ws.Range("E27").Value = '  +0.56%  '
$ws.Range("D28").Value = '2.057'
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("D29").Value = '116.11'
$ws.Range("E29").Value = '  +0.79%  '
$ws.Range("D30").Value = '5.080'
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").Value = '0.08829'
$ws.Range("E31").Value = '  -0.42%  '
$ws.Range("D32").Value = '0.7733'
$ws.Range("E32").Value = '  +5.77%  '
$ws.Range("D33").Value = '3.084'
$ws.Range("E33").Value = '  +4.24%  '
$ws.Range("D34").Value = '1.175'
$ws.Range("E34").Value = '  +3.69%  '
$ws.Range("E35").Value = '  +1.61%  '
$ws.Range("D36").Value = '2.733'
$ws.Range("E36").Value = '  +12.49%  '
$ws.Range("D37").Value = '1.079'
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("D38").Value = '0.01959'
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("D40").Value = '2.967'
$ws.Range("E40").Value = '  +0.81%  '
$ws.Range("D41").Value = '7.043'
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("D42").Value = '0.5141'
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").Value = '0.1643'
$ws.Range("E43").Value = '  +0.85%  '
$ws.Range("D44").Value = '8.428'
$ws.Range("E44").Value = '  +2.55%  '
$ws.Range("D45").Value = '0.4805'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").Value = '10.39'
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("D48").Value = '102.57'
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("D49").Value = '1.647'
$ws.Range("E49").Value = '  +1.23%  '
$ws.Range("D50").Value = '0.06221'
$ws.Range("E50").Value = '  +0.17%  '
$ws.Range("D51").Value = '65.75'
$ws.Range("E51").Value = '  +2.02%  '

$dRange.NumberFormat = $originalFormat
$dRange.Style = "Normal"
